$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns with the refreshed
# crypto-ticker snapshot. Values are entered with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of auto-converting numeric-looking text into numbers/percentages;
# the Style reset afterwards clears the quote-prefix formatting Excel adds
# so the cell keeps its original (default) appearance.
$ws.Range("D2").Value = "'279.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.60%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.77%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.811"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.27%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06351"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.46%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.940"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.56%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.368"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6.36%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8827"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.73%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9555"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'5.10%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'5.68%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.05254"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.94%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.16%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.87%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09070"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.39%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001558"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.14%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0006268"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.82%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005808"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-3.33%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.464"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.40%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'7.03%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3127"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.77%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1340"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.06%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.861"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-6.28%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04317"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.61%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001183"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.44%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'5.46%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-12.80%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04088"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.39%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006732"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'60.98%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1162"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.41%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002289"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'11.04%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01251"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-10.12%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005221"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.92%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'819.38%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.02250"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'6.05%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("E50").Style = "Normal"
